# "added basic data type error"
# Add a new "Int_test" column (G) next to the existing ID/First_Name/Last_Name
# table. The first two data rows hold text values ("abc", "dfe") while the
# third data row holds a genuine number (123) - the intentional "basic data
# type error" referenced in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = "Int_test"
$ws.Range("G5").Value = "abc"
$ws.Range("G6").Value = "dfe"
$ws.Range("G7").Value = 123

$ws.Range("G8").Select() | Out-Null
